$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.057186661957012
$ws.Cells.Item(2, 4).Value = 1.064962594501575
$ws.Cells.Item(2, 5).Value = 1.066467762076237
$ws.Cells.Item(2, 6).Value = 1.07622270676888
$ws.Cells.Item(2, 9).Value = 1.055395891465743
$ws.Cells.Item(2, 10).Value = 1.062184394800426
$ws.Cells.Item(2, 11).Value = 1.067677116467527
$ws.Cells.Item(2, 12).Value = 1.069178231209384
$ws.Cells.Item(2, 13).Value = 1.078907206625973
$ws.Cells.Item(2, 14).Value = 1.024418664857259

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.058222316370262
$ws.Cells.Item(3, 4).Value = 1.065814180466582
$ws.Cells.Item(3, 5).Value = 1.06745584061539
$ws.Cells.Item(3, 6).Value = 1.077209933029698
$ws.Cells.Item(3, 9).Value = 1.055749173331715
$ws.Cells.Item(3, 10).Value = 1.062872005342826
$ws.Cells.Item(3, 11).Value = 1.068343918397649
$ws.Cells.Item(3, 12).Value = 1.069981477550918
$ws.Cells.Item(3, 13).Value = 1.079711480788661
$ws.Cells.Item(3, 14).Value = 1.024654500014195

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.0588925822678
$ws.Cells.Item(4, 4).Value = 1.066365348098445
$ws.Cells.Item(4, 5).Value = 1.068096085510655
$ws.Cells.Item(4, 6).Value = 1.077849352035418
$ws.Cells.Item(4, 9).Value = 1.055976592070734
$ws.Cells.Item(4, 10).Value = 1.063316446477383
$ws.Cells.Item(4, 11).Value = 1.068774871913746
$ws.Cells.Item(4, 12).Value = 1.070501503341489
$ws.Cells.Item(4, 13).Value = 1.080231895723413
$ws.Cells.Item(4, 14).Value = 1.024806785839723

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.059174392597579
$ws.Cells.Item(5, 4).Value = 1.066597090490866
$ws.Cells.Item(5, 5).Value = 1.068365457223547
$ws.Cells.Item(5, 6).Value = 1.078118310909836
$ws.Cells.Item(5, 9).Value = 1.056071916571387
$ws.Cells.Item(5, 10).Value = 1.063503172105453
$ws.Cells.Item(5, 11).Value = 1.068955921747708
$ws.Cells.Item(5, 12).Value = 1.070720186585602
$ws.Cells.Item(5, 13).Value = 1.080450676554002
$ws.Cells.Item(5, 14).Value = 1.024870731082614

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.059221711578859
$ws.Cells.Item(6, 4).Value = 1.066636002886417
$ws.Cells.Item(6, 5).Value = 1.068410698372823
$ws.Cells.Item(6, 6).Value = 1.078163478894596
$ws.Cells.Item(6, 9).Value = 1.056087905410558
$ws.Cells.Item(6, 10).Value = 1.063534517262922
$ws.Cells.Item(6, 11).Value = 1.068986313583692
$ws.Cells.Item(6, 12).Value = 1.070756908222466
$ws.Cells.Item(6, 13).Value = 1.08048741070659
$ws.Cells.Item(6, 14).Value = 1.02488146332325

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.058896347711372
$ws.Cells.Item(7, 4).Value = 1.066368444528652
$ws.Cells.Item(7, 5).Value = 1.068099684033653
$ws.Cells.Item(7, 6).Value = 1.077852945299873
$ws.Cells.Item(7, 9).Value = 1.05597786691109
$ws.Cells.Item(7, 10).Value = 1.063318941976245
$ws.Cells.Item(7, 11).Value = 1.068777291594472
$ws.Cells.Item(7, 12).Value = 1.070504425145445
$ws.Cells.Item(7, 13).Value = 1.080234819090638
$ws.Cells.Item(7, 14).Value = 1.024807640577002

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.057536639568082
$ws.Cells.Item(8, 4).Value = 1.065250363683681
$ws.Cells.Item(8, 5).Value = 1.066801502672676
$ws.Cells.Item(8, 6).Value = 1.076556216187852
$ws.Cells.Item(8, 9).Value = 1.055515528533467
$ws.Cells.Item(8, 10).Value = 1.062416876990828
$ws.Cells.Item(8, 11).Value = 1.0679025710157
$ws.Cells.Item(8, 12).Value = 1.069449635427753
$ws.Cells.Item(8, 13).Value = 1.079179015361888
$ws.Cells.Item(8, 14).Value = 1.024498431593208

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.055141648764631
$ws.Cells.Item(9, 4).Value = 1.063281223512014
$ws.Cells.Item(9, 5).Value = 1.064520807722249
$ws.Cells.Item(9, 6).Value = 1.074275975955999
$ws.Cells.Item(9, 9).Value = 1.05469181667244
$ws.Cells.Item(9, 10).Value = 1.060823596089238
$ws.Cells.Item(9, 11).Value = 1.066357301420617
$ws.Cells.Item(9, 12).Value = 1.067593070299275
$ws.Cells.Item(9, 13).Value = 1.077318548021915
$ws.Cells.Item(9, 14).Value = 1.023951161504689

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.053545659002142
$ws.Cells.Item(10, 4).Value = 1.061969216202544
$ws.Cells.Item(10, 5).Value = 1.063005006260479
$ws.Cells.Item(10, 6).Value = 1.072759058963983
$ws.Cells.Item(10, 9).Value = 1.054136630176501
$ws.Cells.Item(10, 10).Value = 1.059758924038276
$ws.Cells.Item(10, 11).Value = 1.065324525467293
$ws.Cells.Item(10, 12).Value = 1.066356814873566
$ws.Cells.Item(10, 13).Value = 1.076078265765549
$ws.Cells.Item(10, 14).Value = 1.023584711700146

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.052854737291486
$ws.Cells.Item(11, 4).Value = 1.061401288835975
$ws.Cells.Item(11, 5).Value = 1.062349762223532
$ws.Cells.Item(11, 6).Value = 1.072102994515751
$ws.Cells.Item(11, 9).Value = 1.053894797068029
$ws.Cells.Item(11, 10).Value = 1.059297323802264
$ws.Cells.Item(11, 11).Value = 1.064876711234831
$ws.Cells.Item(11, 12).Value = 1.065821854883367
$ws.Cells.Item(11, 13).Value = 1.07554122332027
$ws.Cells.Item(11, 14).Value = 1.02342565731709

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.052598121017824
$ws.Cells.Item(12, 4).Value = 1.061190362989065
$ws.Cells.Item(12, 5).Value = 1.062106542544203
$ws.Cells.Item(12, 6).Value = 1.071859419329528
$ws.Cells.Item(12, 9).Value = 1.053804754237534
$ws.Cells.Item(12, 10).Value = 1.059125776482357
$ws.Cells.Item(12, 11).Value = 1.064710280953251
$ws.Cells.Item(12, 12).Value = 1.065623199544246
$ws.Cells.Item(12, 13).Value = 1.075341743493795
$ws.Cells.Item(12, 14).Value = 1.023366520654153

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.05265316504856
$ws.Cells.Item(13, 4).Value = 1.0612356060385
$ws.Cells.Item(13, 5).Value = 1.062158706417307
$ws.Cells.Item(13, 6).Value = 1.071911661758822
$ws.Cells.Item(13, 9).Value = 1.053824078468256
$ws.Cells.Item(13, 10).Value = 1.059162577963114
$ws.Cells.Item(13, 11).Value = 1.06474598497776
$ws.Cells.Item(13, 12).Value = 1.065665809399245
$ws.Cells.Item(13, 13).Value = 1.075384532511767
$ws.Cells.Item(13, 14).Value = 1.023379208235642

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.052833524845329
$ws.Cells.Item(14, 4).Value = 1.061383853072134
$ws.Cells.Item(14, 5).Value = 1.062329654195684
$ws.Cells.Item(14, 6).Value = 1.072082858137753
$ws.Cells.Item(14, 9).Value = 1.053887358489656
$ws.Cells.Item(14, 10).Value = 1.059283145450066
$ws.Cells.Item(14, 11).Value = 1.064862955933512
$ws.Cells.Item(14, 12).Value = 1.065805432894593
$ws.Cells.Item(14, 13).Value = 1.075524734219299
$ws.Cells.Item(14, 14).Value = 1.023420770219357

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.052944653489329
$ws.Cells.Item(15, 4).Value = 1.061475196644621
$ws.Cells.Item(15, 5).Value = 1.062435002930358
$ws.Cells.Item(15, 6).Value = 1.072188353313502
$ws.Cells.Item(15, 9).Value = 1.053926318881637
$ws.Cells.Item(15, 10).Value = 1.059357419345939
$ws.Cells.Item(15, 11).Value = 1.064935013400622
$ws.Cells.Item(15, 12).Value = 1.065891466526616
$ws.Cells.Item(15, 13).Value = 1.075611117350757
$ws.Cells.Item(15, 14).Value = 1.023446370410537

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.05359151640415
$ws.Cells.Item(16, 4).Value = 1.062006911520582
$ws.Cells.Item(16, 5).Value = 1.063048516050308
$ws.Cells.Item(16, 6).Value = 1.072802616086606
$ws.Cells.Item(16, 9).Value = 1.054152649649917
$ws.Cells.Item(16, 10).Value = 1.059789546510452
$ws.Cells.Item(16, 11).Value = 1.065354232484329
$ws.Cells.Item(16, 12).Value = 1.066392325756794
$ws.Cells.Item(16, 13).Value = 1.076113907710281
$ws.Cells.Item(16, 14).Value = 1.023595259642387

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.053997317160499
$ws.Cells.Item(17, 4).Value = 1.062340490917674
$ws.Cells.Item(17, 5).Value = 1.0634336537342
$ws.Cells.Item(17, 6).Value = 1.073188133624202
$ws.Cells.Item(17, 9).Value = 1.054294237226056
$ws.Cells.Item(17, 10).Value = 1.060060450611509
$ws.Cells.Item(17, 11).Value = 1.065617033000597
$ws.Cells.Item(17, 12).Value = 1.066706594726938
$ws.Cells.Item(17, 13).Value = 1.076429297269418
$ws.Cells.Item(17, 14).Value = 1.02368855253346

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.054234028584595
$ws.Cells.Item(18, 4).Value = 1.062535079510621
$ws.Cells.Item(18, 5).Value = 1.063658404890206
$ws.Cells.Item(18, 6).Value = 1.073413073795768
$ws.Cells.Item(18, 9).Value = 1.054376684523071
$ws.Cells.Item(18, 10).Value = 1.060218407427643
$ws.Cells.Item(18, 11).Value = 1.065770260693621
$ws.Cells.Item(18, 12).Value = 1.066889935923052
$ws.Cells.Item(18, 13).Value = 1.076613259400396
$ws.Cells.Item(18, 14).Value = 1.023742932098589

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.054314743560406
$ws.Cells.Item(19, 4).Value = 1.062601432122385
$ws.Cells.Item(19, 5).Value = 1.063735057331548
$ws.Cells.Item(19, 6).Value = 1.073489785168545
$ws.Cells.Item(19, 9).Value = 1.054404773451232
$ws.Cells.Item(19, 10).Value = 1.060272256951053
$ws.Cells.Item(19, 11).Value = 1.065822497291746
$ws.Cells.Item(19, 12).Value = 1.066952456186196
$ws.Cells.Item(19, 13).Value = 1.076675985844289
$ws.Cells.Item(19, 14).Value = 1.023761467911347

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.053953777058914
$ws.Cells.Item(20, 4).Value = 1.062304699204162
$ws.Cells.Item(20, 5).Value = 1.063392321047617
$ws.Cells.Item(20, 6).Value = 1.07314676354995
$ws.Cells.Item(20, 9).Value = 1.054279060527208
$ws.Cells.Item(20, 10).Value = 1.060031391055441
$ws.Cells.Item(20, 11).Value = 1.065588843134127
$ws.Cells.Item(20, 12).Value = 1.066672873179743
$ws.Cells.Item(20, 13).Value = 1.076395458874952
$ws.Cells.Item(20, 14).Value = 1.023678546873484

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.052780412751413
$ws.Cells.Item(21, 4).Value = 1.061340197247585
$ws.Cells.Item(21, 5).Value = 1.062279309716723
$ws.Cells.Item(21, 6).Value = 1.072032441854277
$ws.Cells.Item(21, 9).Value = 1.053868730041206
$ws.Cells.Item(21, 10).Value = 1.059247643769702
$ws.Cells.Item(21, 11).Value = 1.064828513443275
$ws.Cells.Item(21, 12).Value = 1.06576431579845
$ws.Cells.Item(21, 13).Value = 1.075483448261672
$ws.Cells.Item(21, 14).Value = 1.023408532817044

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.052042803992532
$ws.Cells.Item(22, 4).Value = 1.060733936916183
$ws.Cells.Item(22, 5).Value = 1.061580482616356
$ws.Cells.Item(22, 6).Value = 1.071332497097306
$ws.Cells.Item(22, 9).Value = 1.053609493705616
$ws.Cells.Item(22, 10).Value = 1.058754359178024
$ws.Cells.Item(22, 11).Value = 1.064349931076936
$ws.Cells.Item(22, 12).Value = 1.065193374417861
$ws.Cells.Item(22, 13).Value = 1.074910040761706
$ws.Cells.Item(22, 14).Value = 1.023238435826264

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.052433811870168
$ws.Cells.Item(23, 4).Value = 1.06105531160656
$ws.Cells.Item(23, 5).Value = 1.061950852172604
$ws.Cells.Item(23, 6).Value = 1.07170348701879
$ws.Cells.Item(23, 9).Value = 1.053747037778842
$ws.Cells.Item(23, 10).Value = 1.059015907085379
$ws.Cells.Item(23, 11).Value = 1.064603687049411
$ws.Cells.Item(23, 12).Value = 1.065496012188599
$ws.Cells.Item(23, 13).Value = 1.07521401384832
$ws.Cells.Item(23, 14).Value = 1.023328638528702

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.053973450911065
$ws.Cells.Item(24, 4).Value = 1.06232087188516
$ws.Cells.Item(24, 5).Value = 1.063410997179359
$ws.Cells.Item(24, 6).Value = 1.073165456676113
$ws.Cells.Item(24, 9).Value = 1.054285918651549
$ws.Cells.Item(24, 10).Value = 1.060044521993917
$ws.Cells.Item(24, 11).Value = 1.065601581104862
$ws.Cells.Item(24, 12).Value = 1.06668811039211
$ws.Cells.Item(24, 13).Value = 1.076410748986546
$ws.Cells.Item(24, 14).Value = 1.023683068113622

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.055760693523403
$ws.Cells.Item(25, 4).Value = 1.063790163642956
$ws.Cells.Item(25, 5).Value = 1.06510960367189
$ws.Cells.Item(25, 6).Value = 1.0748649041181
$ws.Cells.Item(25, 9).Value = 1.054905832504664
$ws.Cells.Item(25, 10).Value = 1.055469181667244
$ws.Cells.Item(25, 11).Value = 1.066757250164677
$ws.Cells.Item(25, 12).Value = 1.068072783180826
$ws.Cells.Item(25, 13).Value = 1.077799520720382
$ws.Cells.Item(25, 14).Value = 1.02409292748266
